$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Email column (B2:B6)
$ws.Range("B2").Value = "groblaiaaa@czamm.hmB"
$ws.Range("B3").Value = "legemnimma@crnao.sdl"
$ws.Range("B4").Value = "lgmmma@oeiin.cag"
$ws.Range("B5").Value = "oncgeig@dwlao.imm"
$ws.Range("B6").Value = "trasylgmejai@cmaro.iam"

# Phone column (C2:C6)
$ws.Range("C2").Value = "34-092-5500200"
$ws.Range("C3").Value = "20-029-0503504"
$ws.Range("C4").Value = "55-003-9400022"
$ws.Range("C5").Value = "50-050-4200923"
$ws.Range("C6").Value = "92-040-5253000"

# Address column (D2:D6)
$ws.Range("D2").Value = "srtd/ermtend/seyap"
$ws.Range("D3").Value = "nsre/paedtte/sdyrm"
$ws.Range("D4").Value = "setr/nemtrds/aypde"
$ws.Range("D5").Value = "sade/rtmsyre/dpten"
$ws.Range("D6").Value = "tsdt/yerrpae/dmnse"
